# Applies the "fix c algorithm and add functions to python" edit to the
# lieux.xlsx-style workbook. Rewrites Feuil1 (sheet1) to hold x/y sample
# points and Feuil2 (sheet2) to add a "vitesse" column with two data rows.

$wb = $excel.ActiveWorkbook

# --- Feuil1 : x / y table -------------------------------------------------
$ws1 = $wb.Worksheets.Item("Feuil1")

# Clear out the old A1:K11 block, we only need A:B now.
$ws1.Range("A1:K11").Clear()

$sheet1Data = @(
    @("x", "y"),
    @(0, 0),
    @(-10, -10),
    @(9, 4),
    @(1, 5),
    @(7, 6),
    @(4, 1),
    @(6, 9),
    @(3, 5),
    @(10, 10),
    @(2, 6)
)

for ($i = 0; $i -lt $sheet1Data.Length; $i++) {
    $row = $i + 1
    $ws1.Cells.Item($row, 1).Value = $sheet1Data[$i][0]
    $ws1.Cells.Item($row, 2).Value = $sheet1Data[$i][1]
}

$ws1.Range("B2").Select() | Out-Null

# --- Feuil2 : add "vitesse" column ----------------------------------------
$ws2 = $wb.Worksheets.Item("Feuil2")

$ws2.Cells.Item(1, 4).Value = "vitesse"

$ws2.Cells.Item(2, 1).Value = 10
$ws2.Cells.Item(2, 2).Value = 100
$ws2.Cells.Item(2, 3).Value = 3
$ws2.Cells.Item(2, 4).Value = 50

$ws2.Cells.Item(3, 1).Value = 5
$ws2.Cells.Item(3, 2).Value = 80
$ws2.Cells.Item(3, 3).Value = 2
$ws2.Cells.Item(3, 4).Value = 50

$ws2.Range("D4").Select() | Out-Null
